# Apply the changes described by the diff to the "Overview" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Updated "report date" labels (shared strings used by row 9) ---
# I9 and K9 both showed "1402-02-23 (5)" -> now "1402-03-09 (6)"
$ws.Range("I9").Value = "1402-03-09 (6)"
$ws.Range("K9").Value = "1402-03-09 (6)"
# M9 showed "1402-02-23 (2)" -> now "1402-03-09 (3)"
$ws.Range("M9").Value = "1402-03-09 (3)"

# --- Updated numeric figures in the financial table ---
$ws.Range("M12").Value = -7688060
$ws.Range("M13").Value = 3835994
$ws.Range("I14").Value = -41717
$ws.Range("M14").Value = -315981
$ws.Range("I17").Value = 2695677
$ws.Range("M17").Value = 3986005
$ws.Range("M18").Value = -424091
$ws.Range("I20").Value = 1916955
$ws.Range("M20").Value = 3757801
$ws.Range("M21").Value = -812907
$ws.Range("I22").Value = 1647869
$ws.Range("M22").Value = 2944894
$ws.Range("I24").Value = 1647869
$ws.Range("M24").Value = 2944894
$ws.Range("M25").Value = 206
$ws.Range("M27").Value = 206
